$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Reports")

# Update the test case name for the second Register bug (row 4, column C)
# from the placeholder "Register" to the proper test case id "Register_TC05"
$ws.Range("C4").Value = "Register_TC05"

# Reflect the final cell selection left by the author when the file was saved
$ws.Range("F5").Select()
